# LOB1003.xlsx update
# - Insert a new row for "Docentes responsaveis" data (row 13)
# - Refresh several text cells with the real (non-placeholder) content
# - Tidy up column A's <col> definition (split off from column B)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a new row at 13 (pushes old rows 13-21 down to 14-22) and
#    make it look like the rest of the B/C data rows (no A cell, B/C
#    styled/wrapped the same way as the other description rows).
# ---------------------------------------------------------------------
$ws.Rows("13:13").Insert()
$ws.Range("A13").Clear()

$ws.Range("B11").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Content updates
# ---------------------------------------------------------------------

# Objetivos (row 10) - was showing the teacher's name by mistake
$objPt = 'Fornecer fundamentos teóricos de limite e derivadas, destacando aspectos geométricos e interpretações físicas, elementos fundamentais para estudos de Engenharia'
$ws.Range("B10").Value = $objPt
$ws.Range("C10").Value = $objPt

# Docentes responsaveis (new row 13)
$docente = '5840692 - Diovana Aparecida dos Santos Napoleão'
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente

# Programa resumido (row 14) - was showing "Semestral" by mistake
$programaResumidoPt = 'Números Reais, funções de variável real, limites e derivadas de funções Reais. Aplicações da derivada e Fórmula de Taylor.'
$ws.Range("B14").Value = $programaResumidoPt
$ws.Range("C14").Value = $programaResumidoPt

# Programa (row 16) - was showing a date by mistake
$programaPt = '•Números e Funções Reais: função trigonométrica, exponencial e logarítmica. Função composta e inversa.•Limite: Definição, propriedades algébricas e Teorema do confronto. Limites infinitos e ao infinito.•Continuidade de funções Reais: Teorema de Weierstrass e teorema do valor intermediário.•Derivada de funções Reais: Definição, Interpretação física e geométrica, regras de derivação, regra da cadeia, derivada da função inversa e derivação implícita, Regra de l’ hopital, Teorema do valor Médio e consequências, Formula de Taylor, taxas de variação, máximos e mínimos (otimização).'
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt

# Metodo (row 19)
$metodoTxt = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("B19").Value = $metodoTxt
$ws.Range("C19").Value = $metodoTxt

# Criterio (row 20)
$criterioTxt = 'NF≥ 5,0.'
$ws.Range("B20").Value = $criterioTxt
$ws.Range("C20").Value = $criterioTxt

# Norma de recuperacao (row 21)
$normaTxt = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("B21").Value = $normaTxt
$ws.Range("C21").Value = $normaTxt

# Bibliografia (row 22)
$bibliografia = "STEWART, James. Cálculo São Paulo: Cengage Learning, 2009. v.1.`nANTON, Howard. Cálculo: um novo horizonte. Porto Alegre: Bookman, 2007.`nTHOMAS, George B. Cálculo São Paulo: Pearson Addison  Wesley, 2009. v.1,`nGUIDORIZZI, Hamilton. Um curso de cálculo. Rio de Janeiro: Livros Técnicos e Científicos, 2001. v.1.`nFLEMMING, Diva M.; GONÇALVES, Mirian B. Cálculo A. São Paulo: Pearson Prentice Hall, 2009."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

# ---------------------------------------------------------------------
# 3) Column A width/definition was sharing a single <col> span with
#    column B (min=1 max=2); give column B an explicit width so it is
#    recorded separately from column A.
# ---------------------------------------------------------------------
$ws.Columns("B:B").ColumnWidth = 60
